# Updated cryptos list on Sat Sep 16 10:58:48 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Prices in column D look numeric ("215.31", "0.506", ...) but must be
    # stored as plain text (as in the source data), not auto-coerced into
    # floating point numbers. Forcing the cell to Text format before the
    # assignment keeps it a string; restoring the "Normal" style afterwards
    # drops the now-unneeded direct formatting so the cell's style index
    # matches the original (unstyled) cells.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Rows 13 and 14 swap places (Polkadot <-> WrappedEther) plus updated prices/volumes.
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "4.20"
$ws.Range("E13").Value = "  +2.09%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.637.55"
$ws.Range("E14").Value = "  -0.55%  "

# Price / volume updates for the remaining rows.
Set-TextValue "D2" "26.642.06"
$ws.Range("E2").Value = "  -0.29%  "

Set-TextValue "D3" "1.642.82"
$ws.Range("E3").Value = "  +0.51%  "

$ws.Range("E4").Value = "  +0.09%  "

Set-TextValue "D5" "215.31"
$ws.Range("E5").Value = "  +0.72%  "

Set-TextValue "D6" "0.506"
$ws.Range("E6").Value = "  +1.23%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("E9").Value = "  +0.78%  "

$ws.Range("E10").Value = "  +0.18%  "

$ws.Range("E11").Value = "  -0.11%  "

Set-TextValue "D12" "1.871.23"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("E15").Value = "  +1.09%  "

Set-TextValue "D16" "65.42"
$ws.Range("E16").Value = "  +2.79%  "

Set-TextValue "D17" "26.684.87"
$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("E18").Value = "  +0.43%  "

Set-TextValue "D19" "217.34"
$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("E21").Value = "  +1.00%  "

$ws.Range("E22").Value = "  +2.48%  "

Set-TextValue "D23" "9.52"
$ws.Range("E23").Value = "  +1.53%  "

Set-TextValue "D24" "2.24"
$ws.Range("E24").Value = "  +13.85%  "

Set-TextValue "D25" "145.48"
$ws.Range("E25").Value = "  -1.27%  "

$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("E28").Value = "  +4.33%  "

$ws.Range("E29").Value = "  +1.28%  "

Set-TextValue "D30" "0.0517"
$ws.Range("E30").Value = "  +2.35%  "

$ws.Range("E31").Value = "  +0.64%  "

Set-TextValue "D32" "3.39"
$ws.Range("E32").Value = "  +2.22%  "

$ws.Range("E33").Value = "  +2.06%  "

Set-TextValue "D34" "1.277.49"
$ws.Range("E34").Value = "  +3.94%  "

$ws.Range("E35").Value = "  +2.46%  "

$ws.Range("E36").Value = "  +4.87%  "

$ws.Range("E37").Value = "  +0.34%  "

Set-TextValue "D38" "0.536"
$ws.Range("E38").Value = "  +6.84%  "

Set-TextValue "D39" "0.829"
$ws.Range("E39").Value = "  +2.60%  "

Set-TextValue "D41" "0.817"
$ws.Range("E41").Value = "  +2.66%  "

$ws.Range("E42").Value = "  -1.36%  "

$ws.Range("E43").Value = "  +1.87%  "

Set-TextValue "D44" "1.781.68"
$ws.Range("E44").Value = "  +0.66%  "

$ws.Range("E45").Value = "  -0.77%  "

Set-TextValue "D46" "59.90"
$ws.Range("E46").Value = "  +8.01%  "

$ws.Range("E47").Value = "  +1.51%  "

$ws.Range("E48").Value = "  +0.57%  "

Set-TextValue "D49" "7.78"
$ws.Range("E49").Value = "  +1.34%  "

Set-TextValue "D50" "0.0970"
$ws.Range("E50").Value = "  +2.82%  "

$ws.Range("E51").Value = "  -0.57%  "
